# Swap the data (columns B:AD) between pairs of rows. The row-number
# column (A) stays put; everything else - id, Div, Date, HomeTeam,
# AwayTeam, scores, odds, etc. - trades places between the two rows
# in each pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

$pairs = @(
    @(73, 74),
    @(124, 125),
    @(159, 160)
)

foreach ($pair in $pairs) {
    $row1 = $pair[0]
    $row2 = $pair[1]

    $vals1 = @{}
    $vals2 = @{}

    foreach ($c in $cols) {
        $vals1[$c] = $ws.Range("$c$row1").Value2
        $vals2[$c] = $ws.Range("$c$row2").Value2
    }

    foreach ($c in $cols) {
        $ws.Range("$c$row1").Value = $vals2[$c]
        $ws.Range("$c$row2").Value = $vals1[$c]
    }
}
